# Update weights and eggs
# Appends two new daily log rows (2021-02-09 and 2021-02-10) to Sheet1,
# matching the existing table's layout (date in col A formatted as a date,
# chinups in col B, pushups in col C), and moves the active selection to
# reflect where the user ended up after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new rows by copying the last existing row (21) so the new
# rows inherit its formatting (in particular column A's date number format)
# instead of creating a brand-new style entry.
$ws.Rows("21").Copy()
$ws.Rows("22").Insert()

$ws.Rows("21").Copy()
$ws.Rows("23").Insert()

# Row 22: 2021-02-09 -> 44236, no chinups, no pushups
$ws.Range("A22").Value = 44236
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0

# Row 23: 2021-02-10 -> 44237, 40 chinups, no pushups
$ws.Range("A23").Value = 44237
$ws.Range("B23").Value = 40
$ws.Range("C23").Value = 0

# Reflect the final cursor position/selection left after the data entry.
$ws.Range("C24").Select()
